# Trade #13 closed at 2026-02-17 20:03:25 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.68   # Current Capital
$summary.Range("B4").Value = -0.33     # Total P&L $
$summary.Range("B5").Value = -0.51     # Total P&L %
$summary.Range("B6").Value = 13        # Total Trades
$summary.Range("B7").Value = 4         # Winning Trades
$summary.Range("B9").Value = 30.77     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.68      # Capital
$status.Range("D5").Value = 13         # Trades
$status.Range("E5").Value = -0.33      # P&L $
$status.Range("F5").Value = -0.32      # P&L %
$status.Range("G5").Value = 30.77      # Win Rate %

# ---------------------------------------------------------------
# New trade row (#13) appended to both "All Trades" and
# "MarketMaking" sheets as row 14.
# ---------------------------------------------------------------
$newRow = @{
    A = 13
    B = "2026-02-17"
    C = "20:03:18"
    D = "MarketMaking"
    E = "UP"
    F = 0.81
    G = 0.83
    H = "CLOSED"
    I = 2.4691
    J = 0.02
    K = 99.68
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(14, 1).Value = $newRow.A

    $ws.Cells.Item(14, 2).NumberFormat = "@"
    $ws.Cells.Item(14, 2).Value = $newRow.B

    $ws.Cells.Item(14, 3).NumberFormat = "@"
    $ws.Cells.Item(14, 3).Value = $newRow.C

    $ws.Cells.Item(14, 4).NumberFormat = "@"
    $ws.Cells.Item(14, 4).Value = $newRow.D

    $ws.Cells.Item(14, 5).NumberFormat = "@"
    $ws.Cells.Item(14, 5).Value = $newRow.E

    $ws.Cells.Item(14, 6).Value = $newRow.F
    $ws.Cells.Item(14, 7).Value = $newRow.G

    $ws.Cells.Item(14, 8).NumberFormat = "@"
    $ws.Cells.Item(14, 8).Value = $newRow.H

    $ws.Cells.Item(14, 9).Value = $newRow.I
    $ws.Cells.Item(14, 10).Value = $newRow.J
    $ws.Cells.Item(14, 11).Value = $newRow.K
    $ws.Cells.Item(14, 12).Value = $newRow.L
    $ws.Cells.Item(14, 13).Value = $newRow.M
    $ws.Cells.Item(14, 14).Value = $newRow.N

    $ws.Cells.Item(14, 15).NumberFormat = "@"
    $ws.Cells.Item(14, 15).Value = $newRow.O

    $ws.Cells.Item(14, 16).NumberFormat = "@"
    $ws.Cells.Item(14, 16).Value = $newRow.P

    $ws.Cells.Item(14, 17).Value = $newRow.Q
}
